$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 753, pushing existing rows 753:835 down to 754:836.
$ws.Rows.Item(753).Insert()

# Populate the newly inserted row 753 with the new weekly record.
$ws.Range("A753").Value = 6
$ws.Range("B753").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C753").Value = "Metropolitana"
$ws.Range("D753").Value = 45212
$ws.Range("E753").Value = 13
$ws.Range("F753").Value = 100112012
$ws.Range("G753").Value = "Espinaca"
$ws.Range("H753").Value = "Sin especificar"
$ws.Range("I753").Value = "Primera"
$ws.Range("J753").Value = 650
$ws.Range("K753").Value = 6000
$ws.Range("L753").Value = 7000
$ws.Range("M753").Value = 6538
$ws.Range("N753").Value = "`$/cuna 10 kilos"
$ws.Range("O753").Value = "Región Metropolitana"
$ws.Range("P753").Value = 654
$ws.Range("Q753").Value = 10
$ws.Range("R753").Value = "Hortaliza"
